# Data augmentation: append new image-to-constellation mapping rows
# (rows 91-113) to Sheet1, extending the used range from A1:B90 to A1:B113.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A91").Value = "image89.jpg"
$ws.Range("B91").Value = "Andromeda"
$ws.Range("A92").Value = "image90.jpg"
$ws.Range("B92").Value = "Andromeda"
$ws.Range("A93").Value = "image91.jpg"
$ws.Range("B93").Value = "Andromeda"
$ws.Range("A94").Value = "image92.jpg"
$ws.Range("B94").Value = "Andromeda"
$ws.Range("A95").Value = "image93.jpg"
$ws.Range("B95").Value = "Antlia"
$ws.Range("A96").Value = "image94.jpg"
$ws.Range("B96").Value = "Antlia"
$ws.Range("A97").Value = "image95.jpg"
$ws.Range("B97").Value = "Antlia"
$ws.Range("A98").Value = "image96.jpg"
$ws.Range("B98").Value = "Antlia"
$ws.Range("A99").Value = "image97.jpg"
$ws.Range("B99").Value = "Antlia"
$ws.Range("A100").Value = "image98.jpg"
$ws.Range("B100").Value = "Apus"
$ws.Range("A101").Value = "image99.jpg"
$ws.Range("B101").Value = "Apus"
$ws.Range("A102").Value = "image100.jpg"
$ws.Range("B102").Value = "Apus"
$ws.Range("A103").Value = "image101.jpg"
$ws.Range("B103").Value = "Apus"
$ws.Range("A104").Value = "image102.jpg"
$ws.Range("B104").Value = "Apus"
$ws.Range("A105").Value = "image103.jpg"
$ws.Range("B105").Value = "Apus"
$ws.Range("A106").Value = "image104.jpg"
$ws.Range("B106").Value = "Aquarius"
$ws.Range("A107").Value = "image105.jpg"
$ws.Range("B107").Value = "Aquarius"
$ws.Range("A108").Value = "image106.jpg"
$ws.Range("B108").Value = "Aquarius"
$ws.Range("A109").Value = "image107.jpg"
$ws.Range("B109").Value = "Aquarius"
$ws.Range("A110").Value = "image108.jpg"
$ws.Range("B110").Value = "Aquarius"
$ws.Range("A111").Value = "image109.jpg"
$ws.Range("B111").Value = "Aquila"
$ws.Range("A112").Value = "image110.jpg"
$ws.Range("B112").Value = "Aquila"
$ws.Range("A113").Value = "image111.jpg"
$ws.Range("B113").Value = "Aquila"
